$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Density (kg/m3)" column (C) header and values, keeping the cell style.
$ws.Range("C1:C6").ClearContents()

# Update the active selection to D10 (per the saved sheet view state in the diff).
$ws.Range("D10").Select()
